# WIP : add experiment tab
#
# The last line of the "update_table" call-out box (the trailing TextBox on
# the final slide) was originally authored as three separate runs:
#
#   [update_table] [(<U+2018>laboratory<U+2019>, ] [connection, grandparent)]
#
# The edit collapses the 2nd and 3rd runs into a single run - the visible
# text does not change, only the run/formatting split moves - keeping the
# formatting ("dirty=0") that the trailing run already carried:
#
#   [update_table] [(<U+2018>laboratory<U+2019>, connection, grandparent)]
#
# NOTE: reading TextRange.Text back through this COM layer transliterates
# the curly quotes to plain ASCII (` and ') even though the underlying
# OOXML keeps the real U+2018/U+2019 glyphs, so matching/searching is done
# against the ASCII form while any text that gets (re)written uses the
# genuine Unicode quote characters so the saved XML keeps its original
# glyphs.

$p = $ppt.ActivePresentation

$openQuoteAscii  = [char]96   # what TextRange.Text reads back for U+2018
$closeQuoteAscii = [char]39   # what TextRange.Text reads back for U+2019
$openQuoteReal   = [char]0x2018
$closeQuoteReal  = [char]0x2019

$marker      = "update_table"
$afterMarker = "(" + $openQuoteAscii + "laboratory" + $closeQuoteAscii + ", connection, grandparent)"
$searchPattern = $marker + $afterMarker

$foundShape = $null
$foundSlideIndex = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                $txt = $shape.TextFrame.TextRange.Text
                if ($txt.Contains($searchPattern)) {
                    $foundShape = $shape
                    $foundSlideIndex = $si
                    break
                }
            }
        }
    }
    if ($foundShape -ne $null) { break }
}

if ($foundShape -eq $null) {
    Write-Host "Could not locate the update_table('laboratory', ...) text box"
} else {
    $tr = $foundShape.TextFrame.TextRange
    $full = $tr.Text

    $markerIdx = $full.IndexOf($marker)
    $run2Start = $markerIdx + $marker.Length

    # Run 2, as originally authored: "(<U+2018>laboratory<U+2019>, "
    $run2TextAscii = "(" + $openQuoteAscii + "laboratory" + $closeQuoteAscii + ", "
    $run2TextReal  = "(" + $openQuoteReal  + "laboratory" + $closeQuoteReal  + ", "
    $run2Len = $run2TextAscii.Length

    # Run 3 follows immediately and already reads "connection, grandparent)".
    $run3Start = $run2Start + $run2Len
    $run3Len   = "connection, grandparent)".Length
    $run3 = $tr.Characters($run3Start + 1, $run3Len)

    # Prepend run 2's text onto run 3 - the newly inserted characters adopt
    # run 3's formatting (dirty="0"), matching the target XML - then delete
    # the now-duplicated, separately-formatted copy of run 2's text that
    # precedes it.
    $run3.InsertBefore($run2TextReal) | Out-Null

    $oldRun2 = $tr.Characters($run2Start + 1, $run2Len)
    $oldRun2.Delete() | Out-Null

    Write-Host "Slide $foundSlideIndex shape '$($foundShape.Name)' updated:"
    Write-Host $tr.Text
}
